# Applies the edits described by the commit:
#  "updated python and expanded livestock elasticities to include all
#   categories to drop category switching"
#
#  1) Sheet "strategy_id-0": lower two blocks of "initial composition"
#     category-count inputs (column C) by 2, rows 107-119 (65 -> 63) and
#     rows 120-132 (67 -> 64).
#  2) Sheets "strategy_id-5006", "strategy_id-5007", "strategy_id-5009":
#     the linear ramp of elasticity-phase-in fractions in columns U:BM
#     (45 columns) is recomputed from a 0.85-anchored ramp to a
#     0.02-step / 0.9-anchored ramp, for the relevant rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) strategy_id-0 : column C, rows 107-132
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("strategy_id-0")

for ($r = 107; $r -le 119; $r++) {
    $wsMain.Range("C$r").Value = 63
}
for ($r = 120; $r -le 132; $r++) {
    $wsMain.Range("C$r").Value = 64
}

# ---------------------------------------------------------------------
# 2) Recomputed ramp values for columns U:BM (45 values), same for every
#    affected row on every affected sheet.
# ---------------------------------------------------------------------
$cols = @("U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH", `
          "AI","AJ","AK","AL","AM","AN","AO","AP","AQ","AR","AS","AT","AU", `
          "AV","AW","AX","AY","AZ","BA","BB","BC","BD","BE","BF","BG","BH", `
          "BI","BJ","BK","BL","BM")

$newVals = @(0.02,0.04,0.06,0.08,0.09999999999999999,0.12,0.14,0.16,0.18, `
             0.2,0.22,0.24,0.26,0.28,0.3,0.32,0.34,0.36,0.38,0.4,0.42,0.44, `
             0.46,0.48,0.5,0.52,0.54,0.5600000000000001,0.5800000000000001, `
             0.6,0.62,0.64,0.6599999999999999,0.68,0.7000000000000001, `
             0.7200000000000001,0.74,0.76,0.78,0.7999999999999999,0.82, `
             0.8400000000000001,0.8600000000000001,0.88,0.9)

$targets = @(
    @{ Sheet = "strategy_id-5006"; Rows = @(6,7,8,9) },
    @{ Sheet = "strategy_id-5007"; Rows = @(41,42,43,44) },
    @{ Sheet = "strategy_id-5009"; Rows = @(41,42,43,44) }
)

foreach ($t in $targets) {
    $ws = $wb.Worksheets.Item($t.Sheet)
    foreach ($r in $t.Rows) {
        for ($i = 0; $i -lt $cols.Length; $i++) {
            $ws.Range("$($cols[$i])$r").Value = $newVals[$i]
        }
    }
}
